$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 277, shifting existing rows 277-331 down to 278-332
$ws.Rows(277).Insert()

# Populate the newly inserted row 277 with the new record's data
$ws.Cells.Item(277, 1).Value = 4
$ws.Cells.Item(277, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(277, 3).Value = "Los Lagos"
$ws.Cells.Item(277, 4).Value = 44711
$ws.Cells.Item(277, 5).Value = 10
$ws.Cells.Item(277, 6).Value = 100112008
$ws.Cells.Item(277, 7).Value = "Coliflor"
$ws.Cells.Item(277, 8).Value = "Sin especificar"
$ws.Cells.Item(277, 9).Value = "Primera"
$ws.Cells.Item(277, 10).Value = 750
$ws.Cells.Item(277, 11).Value = 1500
$ws.Cells.Item(277, 12).Value = 1700
$ws.Cells.Item(277, 13).Value = 1567
$ws.Cells.Item(277, 14).Value = "`$/unidad"
$ws.Cells.Item(277, 15).Value = "Región del Maule"
$ws.Cells.Item(277, 16).Value = 1567
$ws.Cells.Item(277, 17).Value = 1
$ws.Cells.Item(277, 18).Value = "Hortaliza"
